$d = $word.ActiveDocument

$d.Content.Find.Execute("our key contribution is to show that ", $true, $false, $false, $false, $false, $true, 1, $false, "our key contribution is showing that ", 2)

$d.Content.Find.Execute("low validation and generalization error. The major advantage this gradient-based method is", $true, $false, $false, $false, $false, $true, 1, $false, "low validation (and generalization) error. The major advantage of this gradient-based method is", 2)
